$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Row 3 ("Illinois" / "www.illinoisstate.edu") was a stray, misaligned entry that
# knocked the Company Name / Domain Name columns out of step with their hyperlinks.
# Stash the hyperlink cell's pristine formatting before touching anything, so we can
# restore it later (adding a hyperlink via COM redefines the cell's style).
$ws.Range("B2").Copy()
$ws.Range("D1").PasteSpecial($xlPasteFormats)

# Delete the whole row - everything below shifts up one, putting the real data back
# in straight columns (matches the commit message: "have columns straight").
$ws.Rows("3").Delete()

# The row shift does not carry the per-cell hyperlinks along with it, so the
# now-stale/duplicated hyperlink set (including the one for the row that fell off
# the bottom) needs to be rebuilt to match the shifted data in B2:B6.
$ws.Hyperlinks.Delete()

$links = @(
  @{addr = "B2"; url = "http://www.oakstreethealth.com"},
  @{addr = "B3"; url = "http://www.mcmaster.com"},
  @{addr = "B4"; url = "http://www.rewardsnetwork.com"},
  @{addr = "B5"; url = "http://www.gflesch.com"},
  @{addr = "B6"; url = "http://www.sourcebooks.com"}
)
foreach ($link in $links) {
    $ws.Hyperlinks.Add($ws.Range($link.addr), $link.url)
}

# Restore the original hyperlink-cell formatting over the rebuilt range.
$ws.Range("D1").Copy()
$ws.Range("B2:B6").PasteSpecial($xlPasteFormats)
$ws.Range("D1").Clear()
